$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B67: currently stored as text "4" - convert to a real number
$ws.Range("B67").Value = 4

# Add new row 68 with the new annotation data
$ws.Range("A68").Value = "Ying Tang"

# B68 keeps its numeric-looking value stored as text (matches source data)
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "1"

$ws.Range("C68").Value = "really bad"
$ws.Range("D68").Value = "CRT"
$ws.Range("E68").Value = "WRI"
$ws.Range("F68").Value = "3222e19c-371b-4610-a09f-eba8d4490b26"
$ws.Range("G68").Value = "rJTGkKxAZ_annotated.xlsx"
$ws.Range("H68").Value = "This section is really bad."
